{"js": "// Update the date heading and the 25 division problems/answers in the\n// worksheet table. Replacements are applied by position (row/col) since\n// several cell values collide with each other across the edit (e.g. the\n// string that used to be in one cell reappears, unchanged, in a different\n// cell), so a naive global text search-and-replace would be unsafe.\n\n// 1) Update the date paragraph at the top of the document.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateRange = paragraphs.items[0].getRange();\ndateRange.insertText(\"2025-07-08 Tuesday\", Word.InsertLocation.replace);\n\n// 2) Update the practice table. Only rows 0, 4, 8, 12 and 16 hold the 5\n// divison problems each (the remaining rows are blank \"show your work\"\n// rows), giving 25 values total, matching the diff in row-major order.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"57\u00f74=14, 1\", \"95\u00f73=31, 2\", \"69\u00f78=8, 5\", \"60\u00f77=8, 4\", \"70\u00f72=35, 0\"],\n  [\"59\u00f73=19, 2\", \"93\u00f78=11, 5\", \"92\u00f79=10, 2\", \"81\u00f76=13, 3\", \"50\u00f74=12, 2\"],\n  [\"14\u00f76=2, 2\", \"20\u00f79=2, 2\", \"70\u00f77=10, 0\", \"16\u00f78=2, 0\", \"66\u00f79=7, 3\"],\n  [\"85\u00f78=10, 5\", \"10\u00f74=2, 2\", \"45\u00f78=5, 5\", \"99\u00f78=12, 3\", \"45\u00f73=15, 0\"],\n  [\"62\u00f76=10, 2\", \"67\u00f73=22, 1\", \"45\u00f76=7, 3\", \"30\u00f79=3, 3\", \"55\u00f77=7, 6\"],\n];\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\nfor (let i = 0; i < dataRowIndexes.length; i++) {\n  const rowIndex = dataRowIndexes[i];\n  const rowValues = newValues[i];\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const cellPara = cell.body.paragraphs.getFirst();\n    const cellRange = cellPara.getRange();\n    cellRange.insertText(rowValues[col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division problems/answers in the\n# worksheet table. Replacements are applied by position (row/col) since\n# several cell values collide with each other across the edit (e.g. the\n# string that used to be in one cell reappears, unchanged, in a different\n# cell), so a naive global text search-and-replace would be unsafe.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph at the top of the document.\n$d.Paragraphs.Item(1).Range.Text = \"2025-07-08 Tuesday\"\n\n# 2) Update the practice table. Only rows 1, 5, 9, 13 and 17 (1-based) hold\n# the 5 division problems each (the remaining rows are blank \"show your\n# work\" rows), giving 25 values total, matching the diff in row-major order.\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"57\u00f74=14, 1\", \"95\u00f73=31, 2\", \"69\u00f78=8, 5\", \"60\u00f77=8, 4\", \"70\u00f72=35, 0\"),\n    @(\"59\u00f73=19, 2\", \"93\u00f78=11, 5\", \"92\u00f79=10, 2\", \"81\u00f76=13, 3\", \"50\u00f74=12, 2\"),\n    @(\"14\u00f76=2, 2\", \"20\u00f79=2, 2\", \"70\u00f77=10, 0\", \"16\u00f78=2, 0\", \"66\u00f79=7, 3\"),\n    @(\"85\u00f78=10, 5\", \"10\u00f74=2, 2\", \"45\u00f78=5, 5\", \"99\u00f78=12, 3\", \"45\u00f73=15, 0\"),\n    @(\"62\u00f76=10, 2\", \"67\u00f73=22, 1\", \"45\u00f76=7, 3\", \"30\u00f79=3, 3\", \"55\u00f77=7, 6\")\n)\n$dataRowIndexes = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $dataRowIndexes.Length; $i++) {\n    $rowIndex = $dataRowIndexes[$i]\n    $rowValues = $newValues[$i]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $cell = $tbl.Cell($rowIndex, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
